$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D values ("U") for rows 2-6, and clear old F/G values
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 4).Value = "U"          # D column
    $ws.Cells.Item($r, 6).Value = $null        # F column cleared
    $ws.Cells.Item($r, 7).Value = $null        # G column cleared
}

# Update the active selection shown in the sheet view
$ws.Range("D7").Select()
